$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 43199.5621566898

$ws.Range("B3").Value = $newValue
$ws.Range("C3").Value = $newValue
$ws.Range("D3").Value = $newValue
$ws.Range("E3").Value = $newValue
$ws.Range("F3").Value = $newValue
$ws.Range("G3").Value = $newValue
